# Updates market-price-derived columns (H:N) on several Leve profit sheets.
# Values mirror a scheduled market-data refresh; row/column layout is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 32 - Automata for the People
$ws.Range("H32").Value = 15766.667
$ws.Range("I32").Value = 14149
$ws.Range("J32").Value = 19002
$ws.Range("K32").Value = 14149
$ws.Range("L32").Value = 19002
$ws.Range("M32").Value = -13823
$ws.Range("N32").Value = -19654
# row 46 - Always Have an Exit Plan
$ws.Range("H46").Value = 3825
$ws.Range("J46").Value = 4000
$ws.Range("L46").Value = 12000
$ws.Range("N46").Value = -12238
# row 60 - Make Up Your Mind or Else
$ws.Range("H60").Value = 3825
$ws.Range("J60").Value = 4000
$ws.Range("L60").Value = 12000
$ws.Range("N60").Value = -12968
# row 80 - Cleansing the Wicked Humours
$ws.Range("H80").Value = 2423.2122
$ws.Range("I80").Value = 1380.0555
$ws.Range("J80").Value = 3675
$ws.Range("K80").Value = 4140.166499999999
$ws.Range("L80").Value = 11025
$ws.Range("M80").Value = -3142.166499999999
$ws.Range("N80").Value = -13021
# row 83 - Washing Away the Sins (L)
$ws.Range("H83").Value = 2423.2122
$ws.Range("I83").Value = 1380.0555
$ws.Range("J83").Value = 3675
$ws.Range("K83").Value = 12420.4995
$ws.Range("L83").Value = 33075
$ws.Range("M83").Value = -7428.4995
$ws.Range("N83").Value = -43059
# row 97 - Materia Worth
$ws.Range("H97").Value = 1991.6666
$ws.Range("J97").Value = 1991.6666
$ws.Range("L97").Value = 5974.9998
$ws.Range("N97").Value = -6966.9998
# row 112 - Making Ends Meet
$ws.Range("H112").Value = 2011.2307
$ws.Range("J112").Value = 2031.68
$ws.Range("L112").Value = 6095.04
$ws.Range("N112").Value = -8311.040000000001
# row 113 - Amaro Kart
$ws.Range("H113").Value = 6148.6313
$ws.Range("I113").Value = 3863
$ws.Range("J113").Value = 7203.5386
$ws.Range("K113").Value = 3863
$ws.Range("L113").Value = 7203.5386
$ws.Range("M113").Value = -609
$ws.Range("N113").Value = -13711.5386
# row 116 - Growing Up
$ws.Range("H116").Value = 5147.4707
$ws.Range("I116").Value = 3645
$ws.Range("J116").Value = 6837.75
$ws.Range("K116").Value = 3645
$ws.Range("L116").Value = 6837.75
$ws.Range("M116").Value = -203
$ws.Range("N116").Value = -13721.75
# row 132 - Fast-forwarding Flora
$ws.Range("H132").Value = 1133.3667
$ws.Range("I132").Value = 1151.7587
$ws.Range("K132").Value = 3455.2761
$ws.Range("M132").Value = -925.2761
# row 139 - Something Salty and Ceremonial
$ws.Range("H139").Value = 43142.145
$ws.Range("J139").Value = 43142.145
$ws.Range("L139").Value = 43142.145
$ws.Range("N139").Value = -53422.145

$ws = $wb.Worksheets.Item("ARM")
# row 33 - A Leg to Stand On
$ws.Range("H33").Value = 12252.875
$ws.Range("I33").Value = 4505.75
$ws.Range("K33").Value = 4505.75
$ws.Range("M33").Value = -4176.75
# row 45 - Hollow Hallmarks
$ws.Range("H45").Value = 2918.8333
$ws.Range("I45").Value = 2758.1875
$ws.Range("J45").Value = 4204
$ws.Range("K45").Value = 2758.1875
$ws.Range("L45").Value = 4204
$ws.Range("M45").Value = -2381.1875
$ws.Range("N45").Value = -4958
# row 107 - Shielding the Realm
$ws.Range("H107").Value = 96250
$ws.Range("J107").Value = 96250
$ws.Range("L107").Value = 96250
$ws.Range("N107").Value = -103930
# row 135 - Forgiveness for My Shins
$ws.Range("H135").Value = 60175.1
$ws.Range("J135").Value = 60175.1
$ws.Range("L135").Value = 60175.1
$ws.Range("N135").Value = -70315.10000000001

$ws = $wb.Worksheets.Item("BSM")
# row 105 - Ingot to Wing It
$ws.Range("H105").Value = 10405.622
$ws.Range("I105").Value = 8492.518
$ws.Range("K105").Value = 8492.518
$ws.Range("M105").Value = -6745.518
# row 135 - Axes to the Maxes
$ws.Range("H135").Value = 48305
$ws.Range("J135").Value = 48305
$ws.Range("L135").Value = 48305
$ws.Range("N135").Value = -58445

$ws = $wb.Worksheets.Item("CRP")
# row 16 - Raise the Roof
$ws.Range("H16").Value = 1368.3334
$ws.Range("J16").Value = 2344.818
$ws.Range("L16").Value = 2344.818
$ws.Range("N16").Value = -2918.818
# row 99 - O Pine
$ws.Range("H99").Value = 1974.6364
$ws.Range("I99").Value = 1622.25
$ws.Range("J99").Value = 2397.5
$ws.Range("K99").Value = 1622.25
$ws.Range("L99").Value = 2397.5
$ws.Range("M99").Value = -124.25
$ws.Range("N99").Value = -5393.5
# row 107 - Built to Last
$ws.Range("H107").Value = 1701.8422
$ws.Range("J107").Value = 2397.111
$ws.Range("L107").Value = 2397.111
$ws.Range("N107").Value = -6237.111
# row 113 - Patient Patients
$ws.Range("H113").Value = 1368.3334
$ws.Range("J113").Value = 2344.818
$ws.Range("L113").Value = 2344.818
$ws.Range("N113").Value = -6684.818
# row 126 - A Better Conductor
$ws.Range("H126").Value = 1974.6364
$ws.Range("I126").Value = 1622.25
$ws.Range("J126").Value = 2397.5
$ws.Range("K126").Value = 4866.75
$ws.Range("L126").Value = 7192.5
$ws.Range("M126").Value = -2396.75
$ws.Range("N126").Value = -12132.5
# row 132 - Hull Lotta Damage
$ws.Range("H132").Value = 3097
$ws.Range("I132").Value = 2747.2285
$ws.Range("K132").Value = 8241.6855
$ws.Range("M132").Value = -5711.6855
# row 134 - Wood You Be Quiet
$ws.Range("H134").Value = 2143.325
$ws.Range("I134").Value = 1383.9333
$ws.Range("J134").Value = 4421.5
$ws.Range("K134").Value = 4151.7999
$ws.Range("L134").Value = 13264.5
$ws.Range("M134").Value = -1616.7999
$ws.Range("N134").Value = -18334.5

$ws = $wb.Worksheets.Item("CUL")
# row 132 - More Mezcal
$ws.Range("H132").Value = 3336.1667
$ws.Range("I132").Value = 2989.7856
$ws.Range("J132").Value = 3639.25
$ws.Range("K132").Value = 26908.0704
$ws.Range("L132").Value = 32753.25
$ws.Range("M132").Value = -24378.0704
$ws.Range("N132").Value = -37813.25

$ws = $wb.Worksheets.Item("GSM")
# row 128 - To Fight at Her Side
$ws.Range("H128").Value = 60000
$ws.Range("J128").Value = 60000
$ws.Range("L128").Value = 60000
$ws.Range("N128").Value = -69960
# row 136 - Shiny and Good
$ws.Range("H136").Value = 19450.096
$ws.Range("J136").Value = 19450.096
$ws.Range("L136").Value = 58350.288
$ws.Range("N136").Value = -63450.288

$ws = $wb.Worksheets.Item("LTW")
# row 46 - Supply Side Logic
$ws.Range("H46").Value = 6374.357
$ws.Range("J46").Value = 7249.5
$ws.Range("L46").Value = 7249.5
$ws.Range("N46").Value = -7625.5
# row 61 - Spelling Me Softly
$ws.Range("H61").Value = 6661.7144
$ws.Range("I61").Value = 5114.5454
$ws.Range("K61").Value = 5114.5454
$ws.Range("M61").Value = -4912.5454
# row 113 - Peace in Rest
$ws.Range("H113").Value = 6661.7144
$ws.Range("I113").Value = 5114.5454
$ws.Range("K113").Value = 5114.5454
$ws.Range("M113").Value = -2944.5454
# row 136 - Respect for Br'aax
$ws.Range("H136").Value = 3140.4348
$ws.Range("I136").Value = 2435.3901
$ws.Range("J136").Value = 4172.8213
$ws.Range("K136").Value = 7306.1703
$ws.Range("L136").Value = 12518.4639
$ws.Range("M136").Value = -4756.1703
$ws.Range("N136").Value = -17618.4639

$ws = $wb.Worksheets.Item("WVR")
# row 43 - Walk Softly and Carry a Big Halberd
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
# row 64 - Ribbon of Remembrance
$ws.Range("H64").Value = 27735.334
$ws.Range("I64").Value = 27735.334
$ws.Range("K64").Value = 27735.334
$ws.Range("M64").Value = -27487.334
# row 67 - The Road Was a Ribbon of Moonlight (L)
$ws.Range("H67").Value = 27735.334
$ws.Range("I67").Value = 27735.334
$ws.Range("K67").Value = 27735.334
$ws.Range("M67").Value = -26877.334
# row 132 - Comfy Cabins
$ws.Range("H132").Value = 2426.9688
$ws.Range("I132").Value = 1757.0416
$ws.Range("K132").Value = 5271.1248
$ws.Range("M132").Value = -2741.1248
